$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B31:K45").Copy()
$ws.Range("B46").PasteSpecial(-4104)
